$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.280.85'
$ws.Range('E2').Value = '  -1.46%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.491.85'
$ws.Range('E3').Value = '  -4.10%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '199.86'
$ws.Range('E5').Value = '  +2.01%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '549.65'
$ws.Range('E6').Value = '  -4.78%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.486.88'
$ws.Range('E7').Value = '  -4.08%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.605'

$ws.Range('E9').Value = '  -0.12%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.651'
$ws.Range('E10').Value = '  -4.33%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '62.73'
$ws.Range('E11').Value = '  +11.26%  '

$ws.Range('E12').Value = '  -7.46%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000268'
$ws.Range('E13').Value = '  -8.88%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.75'
$ws.Range('E14').Value = '  -3.95%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.045.50'
$ws.Range('E15').Value = '  -4.14%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.493.43'
$ws.Range('E16').Value = '  -4.09%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.124'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '66.913.82'
$ws.Range('E18').Value = '  -1.87%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.20'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.73'
$ws.Range('E20').Value = '  -6.74%  '

$ws.Range('E21').Value = '  -5.79%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '388.67'
$ws.Range('E22').Value = '  -3.77%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.01'
$ws.Range('E23').Value = '  -5.91%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.98'
$ws.Range('E24').Value = '  -6.05%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.06'
$ws.Range('E25').Value = '  -4.77%  '

$ws.Range('E26').Value = '  -0.61%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.12'
$ws.Range('E27').Value = '  -4.11%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.79'
$ws.Range('E28').Value = '  -5.84%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.77'
$ws.Range('E29').Value = '  -4.35%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.87'
$ws.Range('E30').Value = '  -2.79%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '675.55'
$ws.Range('E31').Value = '  -2.62%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.94'
$ws.Range('E32').Value = '  -14.87%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.68'
$ws.Range('E33').Value = '  -4.62%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.76'
$ws.Range('E34').Value = '  -1.65%  '

$ws.Range('E35').Value = '  -7.72%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '38.42'
$ws.Range('E36').Value = '  -10.29%  '

$ws.Range('E37').Value = '  +0.00%  '

$ws.Range('E38').Value = '  -5.63%  '

$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.073.48'
$ws.Range('E39').Value = '  -3.86%  '

$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  -0.05%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.130'
$ws.Range('E41').Value = '  -4.35%  '

$ws.Range('E42').Value = '  -5.21%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0₃0670'
$ws.Range('E43').Value = '  -16.23%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.77'
$ws.Range('E44').Value = '  +6.62%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.50'
$ws.Range('E45').Value = '  -12.69%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.68'
$ws.Range('E46').Value = '  -8.96%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0395'
$ws.Range('E47').Value = '  -6.72%  '

$ws.Range('E48').Value = '  -5.02%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '136.56'
$ws.Range('E49').Value = '  -4.54%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.19'
$ws.Range('E50').Value = '  -7.90%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.86'
$ws.Range('E51').Value = '  -7.48%  '
